$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.707522
$ws.Range("H2").Value = 5.122566
$ws.Range("I2").Value = 0.002012043481081613
$ws.Range("J2").Value = 0.002016814216842583
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 1.497831666666666
$ws.Range("N2").Value = 4.493494999999999
$ws.Range("O2").Value = 0.005149690545628369
$ws.Range("P2").Value = 0.005179071718449603
$ws.Range("Q2").Value = 2.55758052313
$ws.Range("R2").Value = 23.01822470817
$ws.Range("S2").Value = 0.00001036140129191917
$ws.Range("T2").Value = 0.00001044522547181651

$ws.Range("G3").Value = 1.707522
$ws.Range("H3").Value = 5.122566
$ws.Range("I3").Value = 0.002012043481081613
$ws.Range("J3").Value = 0.002016814216842583
$ws.Range("O3").Value = 0.001652912115568146
$ws.Range("P3").Value = 0.001662342681559543
$ws.Range("Q3").Value = 0.8209145376339999
$ws.Range("R3").Value = 7.388230838706
$ws.Range("S3").Value = 0.000003325731046929707
$ws.Range("T3").Value = 0.000003352636353433509

$ws.Range("G4").Value = 1.707522
$ws.Range("H4").Value = 5.122566
$ws.Range("I4").Value = 0.002012043481081613
$ws.Range("J4").Value = 0.002016814216842583
$ws.Range("M4").Value = 134.73733
$ws.Range("N4").Value = 404.21199
$ws.Range("O4").Value = 0.4632400087977464
$ws.Range("P4").Value = 0.4658829898925522
$ws.Range("Q4").Value = 230.06695519626
$ws.Range("R4").Value = 2070.60259676634
$ws.Range("S4").Value = 0.0009320590398776947
$ws.Range("T4").Value = 0.0009395994374004289

$ws.Range("G5").Value = 1.707522
$ws.Range("H5").Value = 5.122566
$ws.Range("I5").Value = 0.002012043481081613
$ws.Range("J5").Value = 0.002016814216842583
$ws.Range("M5").Value = 4.9501725
$ws.Range("N5").Value = 9.900345
$ws.Range("O5").Value = 0.01701917317532092
$ws.Range("P5").Value = 0.01141084986016318
$ws.Range("Q5").Value = 8.452528447545
$ws.Range("R5").Value = 50.71517068527
$ws.Range("S5").Value = 0.00003424331644080351
$ws.Range("T5").Value = 0.0000230135642242333

$ws.Range("G6").Value = 1.707522
$ws.Range("H6").Value = 5.122566
$ws.Range("I6").Value = 0.002012043481081613
$ws.Range("J6").Value = 0.002016814216842583
$ws.Range("M6").Value = 149.1924796666667
$ws.Range("N6").Value = 447.577439
$ws.Range("O6").Value = 0.5129382153657362
$ws.Range("P6").Value = 0.5158647458472754
$ws.Range("Q6").Value = 254.749441265386
$ws.Range("R6").Value = 2292.744971388474
$ws.Range("S6").Value = 0.001032053992424266
$ws.Range("T6").Value = 0.001040403353392671

$ws.Range("I7").Value = 0.000420095518708099
$ws.Range("J7").Value = 0.000421091602904573
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 1.497831666666666
$ws.Range("N7").Value = 4.493494999999999
$ws.Range("O7").Value = 0.005149690545628369
$ws.Range("P7").Value = 0.005179071718449603
$ws.Range("Q7").Value = 0.5339984580872221
$ws.Range("R7").Value = 4.805986122784999
$ws.Range("S7").Value = 0.000002163361920951943
$ws.Range("T7").Value = 0.000002180863611479685

$ws.Range("I8").Value = 0.000420095518708099
$ws.Range("J8").Value = 0.000421091602904573
$ws.Range("O8").Value = 0.001652912115568146
$ws.Range("P8").Value = 0.001662342681559543
$ws.Range("S8").Value = 0.0000006943809725685017
$ws.Range("T8").Value = 0.000000699998544354594

$ws.Range("I9").Value = 0.000420095518708099
$ws.Range("J9").Value = 0.000421091602904573
$ws.Range("M9").Value = 134.73733
$ws.Range("N9").Value = 404.21199
$ws.Range("O9").Value = 0.4632400087977464
$ws.Range("P9").Value = 0.4658829898925522
$ws.Range("Q9").Value = 48.03578938006333
$ws.Range("R9").Value = 432.32210442057
$ws.Range("S9").Value = 0.0001946050517822336
$ws.Range("T9").Value = 0.0001961794149798298

$ws.Range("I10").Value = 0.000420095518708099
$ws.Range("J10").Value = 0.000421091602904573
$ws.Range("M10").Value = 4.9501725
$ws.Range("N10").Value = 9.900345
$ws.Range("O10").Value = 0.01701917317532092
$ws.Range("P10").Value = 0.01141084986016318
$ws.Range("Q10").Value = 1.7648074487225
$ws.Range("R10").Value = 10.588844692335
$ws.Range("S10").Value = 0.000007149678383069405
$ws.Range("T10").Value = 0.000004805013058119535

$ws.Range("I11").Value = 0.000420095518708099
$ws.Range("J11").Value = 0.000421091602904573
$ws.Range("M11").Value = 149.1924796666667
$ws.Range("N11").Value = 447.577439
$ws.Range("O11").Value = 0.5129382153657362
$ws.Range("P11").Value = 0.5158647458472754
$ws.Range("Q11").Value = 53.18925742670856
$ws.Range("R11").Value = 478.703316840377
$ws.Range("S11").Value = 0.0002154830456492755
$ws.Range("T11").Value = 0.0002172263127107894

$ws.Range("G12").Value = 487.9781593333334
$ws.Range("H12").Value = 1463.934478
$ws.Range("I12").Value = 0.5750047580041945
$ws.Range("J12").Value = 0.5763681459167976
$ws.Range("K12").Value = 3.0
$ws.Range("L12").Value = 1.0
$ws.Range("M12").Value = 1.497831666666666
$ws.Range("N12").Value = 4.493494999999999
$ws.Range("O12").Value = 0.005149690545628369
$ws.Range("P12").Value = 0.005179071718449603
$ws.Range("Q12").Value = 730.9091396911789
$ws.Range("R12").Value = 6578.182257220609
$ws.Range("S12").Value = 0.002961096565985529
$ws.Range("T12").Value = 0.002985051963932921

$ws.Range("G13").Value = 487.9781593333334
$ws.Range("H13").Value = 1463.934478
$ws.Range("I13").Value = 0.5750047580041945
$ws.Range("J13").Value = 0.5763681459167976
$ws.Range("O13").Value = 0.001652912115568146
$ws.Range("P13").Value = 0.001662342681559543
$ws.Range("Q13").Value = 234.6021691343443
$ws.Range("R13").Value = 2111.419522209098
$ws.Range("S13").Value = 0.0009504323310144631
$ws.Range("T13").Value = 0.0009581213692488312

$ws.Range("G14").Value = 487.9781593333334
$ws.Range("H14").Value = 1463.934478
$ws.Range("I14").Value = 0.5750047580041945
$ws.Range("J14").Value = 0.5763681459167976
$ws.Range("M14").Value = 134.73733
$ws.Range("N14").Value = 404.21199
$ws.Range("O14").Value = 0.4632400087977464
$ws.Range("P14").Value = 0.4658829898925522
$ws.Range("Q14").Value = 65748.87428688793
$ws.Range("R14").Value = 591739.8685819913
$ws.Range("S14").Value = 0.2663652091566091
$ws.Range("T14").Value = 0.2685201150985445

$ws.Range("G15").Value = 487.9781593333334
$ws.Range("H15").Value = 1463.934478
$ws.Range("I15").Value = 0.5750047580041945
$ws.Range("J15").Value = 0.5763681459167976
$ws.Range("M15").Value = 4.9501725
$ws.Range("N15").Value = 9.900345
$ws.Range("O15").Value = 0.01701917317532092
$ws.Range("P15").Value = 0.01141084986016318
$ws.Range("Q15").Value = 2415.576064932485
$ws.Range("R15").Value = 14493.45638959491
$ws.Range("S15").Value = 0.009786105553106881
$ws.Range("T15").Value = 0.0065768503772372

$ws.Range("G16").Value = 487.9781593333334
$ws.Range("H16").Value = 1463.934478
$ws.Range("I16").Value = 0.5750047580041945
$ws.Range("J16").Value = 0.5763681459167976
$ws.Range("M16").Value = 149.1924796666667
$ws.Range("N16").Value = 447.577439
$ws.Range("O16").Value = 0.5129382153657362
$ws.Range("P16").Value = 0.5158647458472754
$ws.Range("Q16").Value = 72802.67161411578
$ws.Range("R16").Value = 655224.044527042
$ws.Range("S16").Value = 0.2949419143974786
$ws.Range("T16").Value = 0.2973280071078341

$ws.Range("G17").Value = 6.022401
$ws.Range("H17").Value = 12.044802
$ws.Range("I17").Value = 0.007096443075116684
$ws.Range("J17").Value = 0.004742179585905576
$ws.Range("K17").Value = 3.0
$ws.Range("L17").Value = 1.0
$ws.Range("M17").Value = 1.497831666666666
$ws.Range("N17").Value = 4.493494999999999
$ws.Range("O17").Value = 0.005149690545628369
$ws.Range("P17").Value = 0.005179071718449603
$ws.Range("Q17").Value = 9.020542927165
$ws.Range("R17").Value = 54.12325756299
$ws.Range("S17").Value = 0.00003654448581151829
$ws.Range("T17").Value = 0.00002456008817717262

$ws.Range("G18").Value = 6.022401
$ws.Range("H18").Value = 12.044802
$ws.Range("I18").Value = 0.007096443075116684
$ws.Range("J18").Value = 0.004742179585905576
$ws.Range("O18").Value = 0.001652912115568146
$ws.Range("P18").Value = 0.001662342681559543
$ws.Range("Q18").Value = 2.895351586897
$ws.Range("R18").Value = 17.372109521382
$ws.Range("S18").Value = 0.00001172979673630004
$ws.Range("T18").Value = 0.000007883127529271196

$ws.Range("G19").Value = 6.022401
$ws.Range("H19").Value = 12.044802
$ws.Range("I19").Value = 0.007096443075116684
$ws.Range("J19").Value = 0.004742179585905576
$ws.Range("M19").Value = 134.73733
$ws.Range("N19").Value = 404.21199
$ws.Range("O19").Value = 0.4632400087977464
$ws.Range("P19").Value = 0.4658829898925522
$ws.Range("Q19").Value = 811.4422309293301
$ws.Range("R19").Value = 4868.653385575981
$ws.Range("S19").Value = 0.003287356352549759
$ws.Range("T19").Value = 0.002209300804089115

$ws.Range("G20").Value = 6.022401
$ws.Range("H20").Value = 12.044802
$ws.Range("I20").Value = 0.007096443075116684
$ws.Range("J20").Value = 0.004742179585905576
$ws.Range("M20").Value = 4.9501725
$ws.Range("N20").Value = 9.900345
$ws.Range("O20").Value = 0.01701917317532092
$ws.Range("P20").Value = 0.01141084986016318
$ws.Range("Q20").Value = 29.8119238141725
$ws.Range("R20").Value = 119.24769525669
$ws.Range("S20").Value = 0.0001207755936242177
$ws.Range("T20").Value = 0.00005411229926469931

$ws.Range("G21").Value = 6.022401
$ws.Range("H21").Value = 12.044802
$ws.Range("I21").Value = 0.007096443075116684
$ws.Range("J21").Value = 0.004742179585905576
$ws.Range("M21").Value = 149.1924796666667
$ws.Range("N21").Value = 447.577439
$ws.Range("O21").Value = 0.5129382153657362
$ws.Range("P21").Value = 0.5158647458472754
$ws.Range("Q21").Value = 898.4969387370131
$ws.Range("R21").Value = 5390.981632422078
$ws.Range("S21").Value = 0.003640036846394889
$ws.Range("T21").Value = 0.002446323266845317

$ws.Range("G22").Value = 352.5860493333333
$ws.Range("H22").Value = 1057.758148
$ws.Range("I22").Value = 0.415466659920899
$ws.Range("J22").Value = 0.4164517686775497
$ws.Range("K22").Value = 3.0
$ws.Range("L22").Value = 1.0
$ws.Range("M22").Value = 1.497831666666666
$ws.Range("N22").Value = 4.493494999999999
$ws.Range("O22").Value = 0.005149690545628369
$ws.Range("P22").Value = 0.005179071718449603
$ws.Range("Q22").Value = 528.1145499163621
$ws.Range("R22").Value = 4753.030949247259
$ws.Range("S22").Value = 0.002139524730618451
$ws.Range("T22").Value = 0.002156833577256214

$ws.Range("G23").Value = 352.5860493333333
$ws.Range("H23").Value = 1057.758148
$ws.Range("I23").Value = 0.415466659920899
$ws.Range("J23").Value = 0.4164517686775497
$ws.Range("O23").Value = 0.001652912115568146
$ws.Range("P23").Value = 0.001662342681559543
$ws.Range("Q23").Value = 169.5105618930075
$ws.Range("R23").Value = 1525.595057037068
$ws.Range("S23").Value = 0.0006867298757978848
$ws.Range("T23").Value = 0.0006922855498836523

$ws.Range("G24").Value = 352.5860493333333
$ws.Range("H24").Value = 1057.758148
$ws.Range("I24").Value = 0.415466659920899
$ws.Range("J24").Value = 0.4164517686775497
$ws.Range("M24").Value = 134.73733
$ws.Range("N24").Value = 404.21199
$ws.Range("O24").Value = 0.4632400087977464
$ws.Range("P24").Value = 0.4658829898925522
$ws.Range("Q24").Value = 47506.50288242161
$ws.Range("R24").Value = 427558.5259417945
$ws.Range("S24").Value = 0.1924607791969276
$ws.Range("T24").Value = 0.1940177951375384

$ws.Range("G25").Value = 352.5860493333333
$ws.Range("H25").Value = 1057.758148
$ws.Range("I25").Value = 0.415466659920899
$ws.Range("J25").Value = 0.4164517686775497
$ws.Range("M25").Value = 4.9501725
$ws.Range("N25").Value = 9.900345
$ws.Range("O25").Value = 0.01701917317532092
$ws.Range("P25").Value = 0.01141084986016318
$ws.Range("Q25").Value = 1745.36176529351
$ws.Range("R25").Value = 10472.17059176106
$ws.Range("S25").Value = 0.007070899033765943
$ws.Range("T25").Value = 0.004752068606378926

$ws.Range("G26").Value = 352.5860493333333
$ws.Range("H26").Value = 1057.758148
$ws.Range("I26").Value = 0.415466659920899
$ws.Range("J26").Value = 0.4164517686775497
$ws.Range("M26").Value = 149.1924796666667
$ws.Range("N26").Value = 447.577439
$ws.Range("O26").Value = 0.5129382153657362
$ws.Range("P26").Value = 0.5158647458472754
$ws.Range("Q26").Value = 52603.18699591366
$ws.Range("R26").Value = 473428.682963223
$ws.Range("S26").Value = 0.2131087270837892
$ws.Range("T26").Value = 0.2973280071078341

